$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Suecia" / "Rusia" ordering (Suecia was row 22, Rusia was row 23; now Rusia is row 22, Suecia is row 23) ---
$ws.Range("A22").Value = "Rusia"
$ws.Range("A23").Value = "Suecia"

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 09:52"

# --- Row 22 (now Rusia) updated stats ---
$ws.Range("B22").Value = 7497
$ws.Range("C22").Value = 1154
$ws.Range("D22").Value = 494
$ws.Range("E22").Value = 6945
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 58

# --- Row 23 (now Suecia) updated stats ---
$ws.Range("B23").Value = 7206
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 205
$ws.Range("E23").Value = 6524
$ws.Range("F23").Value = 590
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 477

# --- Row 31 (Polonia) ---
$ws.Range("D31").Value = 191
$ws.Range("E31").Value = 4115

# --- Row 61 (Marruecos) ---
$ws.Range("B61").Value = 1141
$ws.Range("C61").Value = 21
$ws.Range("D61").Value = 88
$ws.Range("E61").Value = 970
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 83

# --- Row 65 (Moldavia) ---
$ws.Range("E65").Value = 907
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 21

# --- Row 68 (Armenia) ---
$ws.Range("B68").Value = 853
$ws.Range("C68").Value = 20
$ws.Range("D68").Value = 87
$ws.Range("E68").Value = 758
